$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a "2022-Q4" row above the
#    existing "2022-Q3" row, pushing the old rows down by one.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

# Make room: insert a new row 3 (copy formatting down from row 2 first so the
# blank row created by Insert keeps the same look as the surrounding rows).
$summary.Rows.Item(3).Insert()
$summary.Cells.Item(2,1).Copy()
$summary.Cells.Item(3,1).PasteSpecial(-4122)

# Row 3 now holds what used to be row 2's data ("2022-Q3").
$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(3,2).Value = "2022-Q3"
$summary.Cells.Item(3,3).Value = 4
$summary.Cells.Item(3,4).Value = 0.2

# Row 2 becomes the new "2022-Q4" summary entry.
$summary.Cells.Item(2,2).Value = "2022-Q4"
$summary.Cells.Item(2,3).Value = 7
$summary.Cells.Item(2,4).Value = 0.68

# Row 4 (previously row 3, "2021-Q2") keeps its data; only the running index
# in column A needs to be renumbered.
$summary.Cells.Item(4,1).Value = 2

# ---------------------------------------------------------------------------
# 2. Create the new "2022-Q4" detail sheet, positioned right after "总计"
#    (i.e. right before the existing "2022-Q3" sheet).
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")

# Duplicate the "2022-Q3" sheet (placed immediately before it) so the new
# sheet inherits identical column layout / header formatting, then rename.
$q3.Copy($q3)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# The duplicated sheet only has 4 data rows (rows 2-5); the new data needs 7
# (rows 2-8), so extend column A's formatting down to rows 6-8.
$newSheet.Cells.Item(5,1).Copy()
$newSheet.Cells.Item(6,1).PasteSpecial(-4122)
$newSheet.Cells.Item(7,1).PasteSpecial(-4122)
$newSheet.Cells.Item(8,1).PasteSpecial(-4122)

# Columns B, D-G hold numeric-looking text in the source data (fund codes
# with leading zeros, formatted percentages/amounts); force text formatting
# before writing so they are stored as text rather than auto-converted to
# numbers. Column C (fund name) is never numeric and needs no help; column H
# (ranking) is a genuine number and is left as default/general.
$newSheet.Range("B2:B8").NumberFormat = "@"
$newSheet.Range("D2:G8").NumberFormat = "@"

$data = @(
  @(0, "009774", "财通资管优选回报一年持有期混合", "6.85", "94.99", "4.05", "0.2774", 10),
  @(1, "000800", "华商未来主题混合", "4.12", "74.31", "3.15", "0.1298", 9),
  @(2, "014575", "鑫元清洁能源混合C", "1.28", "94.25", "6.74", "0.0863", 8),
  @(3, "006010", "国融融银灵活配置混合C", "2.12", "65.24", "3.45", "0.0731", 6),
  @(4, "011815", "恒越优势精选混合", "2.64", "92.01", "2.47", "0.0652", 5),
  @(5, "014574", "鑫元清洁能源混合A", "0.66", "94.25", "6.74", "0.0445", 8),
  @(6, "006009", "国融融银灵活配置混合A", "0.19", "65.24", "3.45", "0.0066", 6)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r,1).Value = $row[0]
    $newSheet.Cells.Item($r,2).Value = $row[1]
    $newSheet.Cells.Item($r,3).Value = $row[2]
    $newSheet.Cells.Item($r,4).Value = $row[3]
    $newSheet.Cells.Item($r,5).Value = $row[4]
    $newSheet.Cells.Item($r,6).Value = $row[5]
    $newSheet.Cells.Item($r,7).Value = $row[6]
    $newSheet.Cells.Item($r,8).Value = $row[7]
    $r = $r + 1
}
